$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 81.666664
$ws.Range("I2").Value = 81.666664
$ws.Range("K2").Value = 81.666664
$ws.Range("M2").Value = 31.333336
$ws.Range("H11").Value = 665.7222
$ws.Range("I11").Value = 665.7222
$ws.Range("K11").Value = 665.7222
$ws.Range("M11").Value = -525.7222
$ws.Range("H40").Value = 82698.37
$ws.Range("J40").Value = 3901.1143
$ws.Range("L40").Value = 3901.1143
$ws.Range("N40").Value = -4251.1143
$ws.Range("H86").Value = 666669000
$ws.Range("I86").Value = 500003500
$ws.Range("K86").Value = 500003500
$ws.Range("M86").Value = -500002377
$ws.Range("H89").Value = 666669000
$ws.Range("I89").Value = 500003500
$ws.Range("K89").Value = 2500017500
$ws.Range("M89").Value = -2500011884
$ws.Range("H109").Value = 109000
$ws.Range("J109").Value = 109000
$ws.Range("L109").Value = 109000
$ws.Range("N109").Value = -111774
$ws.Range("H113").Value = 2754.75
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H127").Value = 1058
$ws.Range("I127").Value = 1058
$ws.Range("K127").Value = 3174
$ws.Range("M127").Value = 1786

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2611.1333
$ws.Range("I2").Value = 1583.75
$ws.Range("K2").Value = 1583.75
$ws.Range("M2").Value = -1470.75
$ws.Range("H32").Value = 9438.637000000001
$ws.Range("I32").Value = 7290.75
$ws.Range("J32").Value = 15166.333
$ws.Range("K32").Value = 7290.75
$ws.Range("L32").Value = 15166.333
$ws.Range("M32").Value = -7003.75
$ws.Range("N32").Value = -15740.333
$ws.Range("H116").Value = 2611.1333
$ws.Range("I116").Value = 1583.75
$ws.Range("K116").Value = 1583.75
$ws.Range("M116").Value = 710.25
$ws.Range("H122").Value = 5189.032
$ws.Range("I122").Value = 3752.1765
$ws.Range("K122").Value = 11256.5295
$ws.Range("M122").Value = -8806.529500000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2611.1333
$ws.Range("I3").Value = 1583.75
$ws.Range("K3").Value = 1583.75
$ws.Range("M3").Value = -1469.75
$ws.Range("H99").Value = 3362.45
$ws.Range("I99").Value = 2958.25
$ws.Range("K99").Value = 2958.25
$ws.Range("M99").Value = -1460.25
$ws.Range("H134").Value = 2673
$ws.Range("I134").Value = 2426.5
$ws.Range("J134").Value = 2820.9
$ws.Range("K134").Value = 7279.5
$ws.Range("L134").Value = 8462.700000000001
$ws.Range("M134").Value = -4744.5
$ws.Range("N134").Value = -13532.7

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.666664
$ws.Range("I7").Value = 70.59999999999999
$ws.Range("K7").Value = 70.59999999999999
$ws.Range("M7").Value = 42.40000000000001
$ws.Range("H31").Value = 3713.6562
$ws.Range("J31").Value = 4267
$ws.Range("L31").Value = 4267
$ws.Range("N31").Value = -4857
$ws.Range("H34").Value = 3713.6562
$ws.Range("J34").Value = 4267
$ws.Range("L34").Value = 4267
$ws.Range("N34").Value = -4671
$ws.Range("H122").Value = 1412.25
$ws.Range("I122").Value = 1360
$ws.Range("K122").Value = 4080
$ws.Range("M122").Value = -1630
$ws.Range("H132").Value = 2242.5334
$ws.Range("I132").Value = 3102.125
$ws.Range("K132").Value = 9306.375
$ws.Range("M132").Value = -6776.375
$ws.Range("H134").Value = 4643.4287
$ws.Range("I134").Value = 4137.1816
$ws.Range("K134").Value = 12411.5448
$ws.Range("M134").Value = -9876.5448

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 910.4666999999999
$ws.Range("I34").Value = 109
$ws.Range("J34").Value = 1311.2
$ws.Range("K34").Value = 327
$ws.Range("L34").Value = 3933.6
$ws.Range("M34").Value = -243
$ws.Range("N34").Value = -4101.6
$ws.Range("H40").Value = 3739.6667
$ws.Range("I40").Value = 12
$ws.Range("J40").Value = 7467.3335
$ws.Range("K40").Value = 48
$ws.Range("L40").Value = 29869.334
$ws.Range("M40").Value = 21
$ws.Range("N40").Value = -30007.334
$ws.Range("H55").Value = 4491.5557
$ws.Range("I55").Value = 2744.5
$ws.Range("J55").Value = 4990.7144
$ws.Range("K55").Value = 8233.5
$ws.Range("L55").Value = 14972.1432
$ws.Range("M55").Value = -8056.5
$ws.Range("N55").Value = -15326.1432

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13000
$ws.Range("J46").Value = 22000
$ws.Range("L46").Value = 22000
$ws.Range("N46").Value = -22312
$ws.Range("H70").Value = 4469.3
$ws.Range("I70").Value = 4353.143
$ws.Range("J70").Value = 4493.9395
$ws.Range("K70").Value = 4353.143
$ws.Range("L70").Value = 4493.9395
$ws.Range("M70").Value = -4083.143
$ws.Range("N70").Value = -5033.9395
$ws.Range("H73").Value = 4469.3
$ws.Range("I73").Value = 4353.143
$ws.Range("J73").Value = 4493.9395
$ws.Range("K73").Value = 4353.143
$ws.Range("L73").Value = 4493.9395
$ws.Range("M73").Value = -3417.143
$ws.Range("N73").Value = -6365.9395
$ws.Range("H122").Value = 1932.1666
$ws.Range("I122").Value = 1843
$ws.Range("K122").Value = 5529
$ws.Range("M122").Value = -3079
$ws.Range("H132").Value = 3595.6924
$ws.Range("I132").Value = 3479.125
$ws.Range("K132").Value = 10437.375
$ws.Range("M132").Value = -7907.375

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8679.708000000001
$ws.Range("I7").Value = 9148.263000000001
$ws.Range("K7").Value = 9148.263000000001
$ws.Range("M7").Value = -9036.263000000001
$ws.Range("H16").Value = 1705.9166
$ws.Range("I16").Value = 1705.9166
$ws.Range("K16").Value = 1705.9166
$ws.Range("M16").Value = -1535.9166
$ws.Range("H22").Value = 2602.926
$ws.Range("I22").Value = 2078.3076
$ws.Range("J22").Value = 3090.0715
$ws.Range("K22").Value = 2078.3076
$ws.Range("L22").Value = 3090.0715
$ws.Range("M22").Value = -1783.3076
$ws.Range("N22").Value = -3680.0715
$ws.Range("H27").Value = 2602.926
$ws.Range("I27").Value = 2078.3076
$ws.Range("J27").Value = 3090.0715
$ws.Range("K27").Value = 2078.3076
$ws.Range("L27").Value = 3090.0715
$ws.Range("M27").Value = -1971.3076
$ws.Range("N27").Value = -3304.0715
$ws.Range("H46").Value = 3008.0417
$ws.Range("I46").Value = 678.3333
$ws.Range("K46").Value = 678.3333
$ws.Range("M46").Value = -490.3333
$ws.Range("H55").Value = 393.33334
$ws.Range("J55").Value = 628.3333
$ws.Range("L55").Value = 628.3333
$ws.Range("N55").Value = -974.3333
$ws.Range("H126").Value = 8679.708000000001
$ws.Range("I126").Value = 9148.263000000001
$ws.Range("K126").Value = 27444.789
$ws.Range("M126").Value = -24974.789
$ws.Range("H136").Value = 4456.524
$ws.Range("I136").Value = 3363.4285
$ws.Range("K136").Value = 10090.2855
$ws.Range("M136").Value = -7540.2855

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 211.07692
$ws.Range("I113").Value = 206.83333
$ws.Range("K113").Value = 620.49999
$ws.Range("M113").Value = 1549.50001
$ws.Range("H126").Value = 7843.9287
$ws.Range("I126").Value = 7843.9287
$ws.Range("K126").Value = 23531.7861
$ws.Range("M126").Value = -21061.7861
$ws.Range("H136").Value = 3049.0417
$ws.Range("I136").Value = 2355.9092
$ws.Range("J136").Value = 3635.5386
$ws.Range("K136").Value = 7067.7276
$ws.Range("L136").Value = 10906.6158
$ws.Range("M136").Value = -4517.7276
$ws.Range("N136").Value = -16006.6158
